$wb = $excel.ActiveWorkbook

# Update "展览" (Exhibition) sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1375
$ws1.Range("F5").Value = 9

# Update "全部类型" (All Types) sheet - mirrors the same data
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1375
$ws4.Range("F5").Value = 9
